$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 = row "O", Row 18 = row "P" - continuing the 384-well plate
# pattern (rows A..N already present in rows 3..16).
# Copy the formatting (white-on-gray header style) from the existing
# row-label cell A16 ("N") onto the two new row-label cells.
$ws.Range("A16").Copy()
$ws.Range("A17:A18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A17").Value = "O"
$ws.Range("A18").Value = "P"

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $label = "SB-{0:D3}" -f ($i + 1)
    $ws.Range($cols[$i] + "17").Value = $label
    $ws.Range($cols[$i] + "18").Value = $label
}

# Update the active selection to reflect the new extent of the sheet,
# matching the post-edit workbook state.
$ws.Range("D31").Select() | Out-Null
